$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace the "Two-check" y/y/n row with division-by-N formulas ---
$ws.Range("A2").ClearContents()
$ws.Range("B2:D2").ClearContents()

$ws.Range("B2").Formula = "=4/4"

$ws.Range("C2").NumberFormat = "d-mmm"
$ws.Range("C2").Formula = "=4/3"

$ws.Range("D2").Formula = "=4/2"

$ws.Range("E2").NumberFormat = "d-mmm"
$ws.Range("E2").Formula = "=4/1"

$ws.Range("F2").Formula = "=4/0"

# --- Remove the now obsolete rows (content cleared, rows collapse away) ---
$ws.Rows(3).ClearContents()
$ws.Rows(5).ClearContents()
$ws.Rows(8).ClearContents()
$ws.Rows(9).ClearContents()

# --- Row 10: repurposed header cell for the new checklist ---
$ws.Range("A10").Value = "Voorwaarden om te sluiten: "
$ws.Range("B10").Clear()

# --- New checklist rows 11-14 ---
$ws.Range("A11").Value = "Verkooporder gesloten"
$ws.Range("B11").Value = "A"

$ws.Range("A12").Value = "Geen waarschuwingen "
$ws.Range("B12").Value = "B"

$ws.Range("A13").Value = "Goedkeuring projectleider"
$ws.Range("B13").Value = "C"

$ws.Range("A14").Value = "Geen actiepunten"
$ws.Range("B14").Value = "D"

# --- Update the remembered selection on the sheet view ---
$ws.Range("A22").Select()
